# Update test cases for iems-test (mirrors the authored diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37 (iems-connector-test-mysql-int-2): swap the condition from the
#     now-removed "idV='46'" case to what used to be row 38's condition,
#     and carry over its pageIndex value while clearing the old timeout value.
$ws.Range("C37").Value = "optimalizeType='2'"
$ws.Range("I37").Value = 10
$ws.Range("J37").ClearContents()

# --- Row 38 (iems-connector-test-mysql-int-3): now exercises the
#     CIMSOURCE_TRAINCONFIG table via a runid condition instead of the
#     optimalizeType condition moved up to row 37.
$ws.Range("C38").Value = "runid='20230209101323_2_545a4f6e-a81f-11ed-a352-0242ac120005'"
$ws.Range("F38").Value = "CIMSOURCE_TRAINCONFIG"
$ws.Range("I38").Value = 0

# --- Remove the old row 39 (iems-connector-test-mysql-int-4 /
#     train_cfg_timetrg_YN='0') entirely; rows below shift up by one.
$ws.Rows(39).Delete()

# --- The row that is now 40 (previously 41, iems-connector-test-mysql-bit-1)
#     has its test-id column overwritten with the same text as its
#     description column.
$ws.Range("A40").Value = "good request, data retrieved (no schema check, no condition check)"

# --- Sheet view: drop the frozen-pane scroll offset and move the saved
#     selection.
$sheetView = $ws.Application
$ws.Range("B44").Select()
$ws.Application.ActiveWindow.ScrollColumn = 1
